$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the four rows for accounts 004948033 (GUILHERME), 004208733 (REINALDO),
# 004363260 (LARISSA) and 004364200 (BLOCO) — rows 22 through 25.
$ws.Range("A22:C25").EntireRow.Delete()
